$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.364.56"
$ws.Range("D3").Value = "3.379.99"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'580.27"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'178.52"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("E9").Value = "  +8.58%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "'48.48"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").Value = "'687.05"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "3.922.68"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "69.471.95"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "3.380.24"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'17.86"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'11.27"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'0.909"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "'17.19"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").Value = "'5.36"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'101.46"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "'33.48"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").Value = "'6.94"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'3.79"
$ws.Range("E31").Value = "  +16.02%  "
$ws.Range("D32").Value = "'11.05"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'550.30"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'57.91"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "3.604.26"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "'35.26"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  +8.15%  "
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D44").Value = "'0.0425"
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "'2.66"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("D50").Value = "'129.42"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("E51").Value = "  +0.47%  "
